# Update the worksheet date header and regenerate all 100 addition /
# subtraction problems with new random operands (textual find & replace,
# each old value is unique in the document so a single pass is safe).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-25 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2026-02-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("85+6=", $true, $true, $false, $false, $false, $true, 1, $false, "69+2=", 2) | Out-Null
$d.Content.Find.Execute("95-49=", $true, $true, $false, $false, $false, $true, 1, $false, "88-1=", 2) | Out-Null
$d.Content.Find.Execute("98-96=", $true, $true, $false, $false, $false, $true, 1, $false, "68-37=", 2) | Out-Null
$d.Content.Find.Execute("13+58=", $true, $true, $false, $false, $false, $true, 1, $false, "99-13=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $true, $true, $false, $false, $false, $true, 1, $false, "3+60=", 2) | Out-Null
$d.Content.Find.Execute("58+22=", $true, $true, $false, $false, $false, $true, 1, $false, "63-0=", 2) | Out-Null
$d.Content.Find.Execute("54+5=", $true, $true, $false, $false, $false, $true, 1, $false, "88-27=", 2) | Out-Null
$d.Content.Find.Execute("99-14=", $true, $true, $false, $false, $false, $true, 1, $false, "11+76=", 2) | Out-Null
$d.Content.Find.Execute("9+11=", $true, $true, $false, $false, $false, $true, 1, $false, "90-66=", 2) | Out-Null
$d.Content.Find.Execute("95-87=", $true, $true, $false, $false, $false, $true, 1, $false, "86-21=", 2) | Out-Null
$d.Content.Find.Execute("0+28=", $true, $true, $false, $false, $false, $true, 1, $false, "96-65=", 2) | Out-Null
$d.Content.Find.Execute("55+13=", $true, $true, $false, $false, $false, $true, 1, $false, "61-19=", 2) | Out-Null
$d.Content.Find.Execute("47+24=", $true, $true, $false, $false, $false, $true, 1, $false, "77-24=", 2) | Out-Null
$d.Content.Find.Execute("51+29=", $true, $true, $false, $false, $false, $true, 1, $false, "22+22=", 2) | Out-Null
$d.Content.Find.Execute("11-2=", $true, $true, $false, $false, $false, $true, 1, $false, "4+30=", 2) | Out-Null
$d.Content.Find.Execute("43+44=", $true, $true, $false, $false, $false, $true, 1, $false, "89-7=", 2) | Out-Null
$d.Content.Find.Execute("28+68=", $true, $true, $false, $false, $false, $true, 1, $false, "55-2=", 2) | Out-Null
$d.Content.Find.Execute("79-71=", $true, $true, $false, $false, $false, $true, 1, $false, "87-12=", 2) | Out-Null
$d.Content.Find.Execute("86-11=", $true, $true, $false, $false, $false, $true, 1, $false, "9+66=", 2) | Out-Null
$d.Content.Find.Execute("64-30=", $true, $true, $false, $false, $false, $true, 1, $false, "57-44=", 2) | Out-Null
$d.Content.Find.Execute("28+48=", $true, $true, $false, $false, $false, $true, 1, $false, "35+58=", 2) | Out-Null
$d.Content.Find.Execute("10+48=", $true, $true, $false, $false, $false, $true, 1, $false, "82-51=", 2) | Out-Null
$d.Content.Find.Execute("80-0=", $true, $true, $false, $false, $false, $true, 1, $false, "58-48=", 2) | Out-Null
$d.Content.Find.Execute("89-71=", $true, $true, $false, $false, $false, $true, 1, $false, "83-39=", 2) | Out-Null
$d.Content.Find.Execute("63+6=", $true, $true, $false, $false, $false, $true, 1, $false, "88-8=", 2) | Out-Null
$d.Content.Find.Execute("57-2=", $true, $true, $false, $false, $false, $true, 1, $false, "55-8=", 2) | Out-Null
$d.Content.Find.Execute("63-21=", $true, $true, $false, $false, $false, $true, 1, $false, "30-22=", 2) | Out-Null
$d.Content.Find.Execute("2+67=", $true, $true, $false, $false, $false, $true, 1, $false, "20+2=", 2) | Out-Null
$d.Content.Find.Execute("41+30=", $true, $true, $false, $false, $false, $true, 1, $false, "14-7=", 2) | Out-Null
$d.Content.Find.Execute("73+20=", $true, $true, $false, $false, $false, $true, 1, $false, "95-13=", 2) | Out-Null
$d.Content.Find.Execute("31+37=", $true, $true, $false, $false, $false, $true, 1, $false, "59-35=", 2) | Out-Null
$d.Content.Find.Execute("22+20=", $true, $true, $false, $false, $false, $true, 1, $false, "83-39=", 2) | Out-Null
$d.Content.Find.Execute("71-66=", $true, $true, $false, $false, $false, $true, 1, $false, "1+12=", 2) | Out-Null
$d.Content.Find.Execute("21+25=", $true, $true, $false, $false, $false, $true, 1, $false, "54+29=", 2) | Out-Null
$d.Content.Find.Execute("64-17=", $true, $true, $false, $false, $false, $true, 1, $false, "20+43=", 2) | Out-Null
$d.Content.Find.Execute("22+42=", $true, $true, $false, $false, $false, $true, 1, $false, "83-13=", 2) | Out-Null
$d.Content.Find.Execute("4+73=", $true, $true, $false, $false, $false, $true, 1, $false, "83-81=", 2) | Out-Null
$d.Content.Find.Execute("67-48=", $true, $true, $false, $false, $false, $true, 1, $false, "20+39=", 2) | Out-Null
$d.Content.Find.Execute("78-6=", $true, $true, $false, $false, $false, $true, 1, $false, "1+68=", 2) | Out-Null
$d.Content.Find.Execute("92-16=", $true, $true, $false, $false, $false, $true, 1, $false, "14+30=", 2) | Out-Null
$d.Content.Find.Execute("7+19=", $true, $true, $false, $false, $false, $true, 1, $false, "95-33=", 2) | Out-Null
$d.Content.Find.Execute("18+17=", $true, $true, $false, $false, $false, $true, 1, $false, "54-0=", 2) | Out-Null
$d.Content.Find.Execute("79-65=", $true, $true, $false, $false, $false, $true, 1, $false, "27-10=", 2) | Out-Null
$d.Content.Find.Execute("60-40=", $true, $true, $false, $false, $false, $true, 1, $false, "4+14=", 2) | Out-Null
$d.Content.Find.Execute("82-39=", $true, $true, $false, $false, $false, $true, 1, $false, "44+50=", 2) | Out-Null
$d.Content.Find.Execute("79-52=", $true, $true, $false, $false, $false, $true, 1, $false, "19+45=", 2) | Out-Null
$d.Content.Find.Execute("39+18=", $true, $true, $false, $false, $false, $true, 1, $false, "14+20=", 2) | Out-Null
$d.Content.Find.Execute("25-14=", $true, $true, $false, $false, $false, $true, 1, $false, "42+26=", 2) | Out-Null
$d.Content.Find.Execute("70-36=", $true, $true, $false, $false, $false, $true, 1, $false, "14+9=", 2) | Out-Null
$d.Content.Find.Execute("96-77=", $true, $true, $false, $false, $false, $true, 1, $false, "65-15=", 2) | Out-Null
$d.Content.Find.Execute("29+8=", $true, $true, $false, $false, $false, $true, 1, $false, "99-31=", 2) | Out-Null
$d.Content.Find.Execute("68-17=", $true, $true, $false, $false, $false, $true, 1, $false, "79-78=", 2) | Out-Null
$d.Content.Find.Execute("0+1=", $true, $true, $false, $false, $false, $true, 1, $false, "87-3=", 2) | Out-Null
$d.Content.Find.Execute("49-33=", $true, $true, $false, $false, $false, $true, 1, $false, "78-15=", 2) | Out-Null
$d.Content.Find.Execute("3+0=", $true, $true, $false, $false, $false, $true, 1, $false, "29+10=", 2) | Out-Null
$d.Content.Find.Execute("28-9=", $true, $true, $false, $false, $false, $true, 1, $false, "96-66=", 2) | Out-Null
$d.Content.Find.Execute("22+45=", $true, $true, $false, $false, $false, $true, 1, $false, "6+62=", 2) | Out-Null
$d.Content.Find.Execute("88-14=", $true, $true, $false, $false, $false, $true, 1, $false, "19+75=", 2) | Out-Null
$d.Content.Find.Execute("49+3=", $true, $true, $false, $false, $false, $true, 1, $false, "7+26=", 2) | Out-Null
$d.Content.Find.Execute("81-59=", $true, $true, $false, $false, $false, $true, 1, $false, "60-32=", 2) | Out-Null
$d.Content.Find.Execute("22+18=", $true, $true, $false, $false, $false, $true, 1, $false, "75-8=", 2) | Out-Null
$d.Content.Find.Execute("41+15=", $true, $true, $false, $false, $false, $true, 1, $false, "19+44=", 2) | Out-Null
$d.Content.Find.Execute("90+3=", $true, $true, $false, $false, $false, $true, 1, $false, "59+15=", 2) | Out-Null
$d.Content.Find.Execute("45+38=", $true, $true, $false, $false, $false, $true, 1, $false, "91-48=", 2) | Out-Null
$d.Content.Find.Execute("86-18=", $true, $true, $false, $false, $false, $true, 1, $false, "32+34=", 2) | Out-Null
$d.Content.Find.Execute("87-35=", $true, $true, $false, $false, $false, $true, 1, $false, "4+0=", 2) | Out-Null
$d.Content.Find.Execute("89-72=", $true, $true, $false, $false, $false, $true, 1, $false, "31+68=", 2) | Out-Null
$d.Content.Find.Execute("71-29=", $true, $true, $false, $false, $false, $true, 1, $false, "98-40=", 2) | Out-Null
$d.Content.Find.Execute("58-1=", $true, $true, $false, $false, $false, $true, 1, $false, "73-11=", 2) | Out-Null
$d.Content.Find.Execute("3+75=", $true, $true, $false, $false, $false, $true, 1, $false, "30-23=", 2) | Out-Null
$d.Content.Find.Execute("45+9=", $true, $true, $false, $false, $false, $true, 1, $false, "57-41=", 2) | Out-Null
$d.Content.Find.Execute("7+14=", $true, $true, $false, $false, $false, $true, 1, $false, "72+19=", 2) | Out-Null
$d.Content.Find.Execute("42+14=", $true, $true, $false, $false, $false, $true, 1, $false, "29-25=", 2) | Out-Null
$d.Content.Find.Execute("14+49=", $true, $true, $false, $false, $false, $true, 1, $false, "93+3=", 2) | Out-Null
$d.Content.Find.Execute("69-5=", $true, $true, $false, $false, $false, $true, 1, $false, "8+53=", 2) | Out-Null
$d.Content.Find.Execute("74-44=", $true, $true, $false, $false, $false, $true, 1, $false, "86-61=", 2) | Out-Null
$d.Content.Find.Execute("35-8=", $true, $true, $false, $false, $false, $true, 1, $false, "2+9=", 2) | Out-Null
$d.Content.Find.Execute("84-80=", $true, $true, $false, $false, $false, $true, 1, $false, "93-50=", 2) | Out-Null
$d.Content.Find.Execute("57-36=", $true, $true, $false, $false, $false, $true, 1, $false, "86-28=", 2) | Out-Null
$d.Content.Find.Execute("12+34=", $true, $true, $false, $false, $false, $true, 1, $false, "12+7=", 2) | Out-Null
$d.Content.Find.Execute("5+60=", $true, $true, $false, $false, $false, $true, 1, $false, "16+46=", 2) | Out-Null
$d.Content.Find.Execute("19-11=", $true, $true, $false, $false, $false, $true, 1, $false, "76-50=", 2) | Out-Null
$d.Content.Find.Execute("16+55=", $true, $true, $false, $false, $false, $true, 1, $false, "81-41=", 2) | Out-Null
$d.Content.Find.Execute("45+30=", $true, $true, $false, $false, $false, $true, 1, $false, "3+78=", 2) | Out-Null
$d.Content.Find.Execute("37+39=", $true, $true, $false, $false, $false, $true, 1, $false, "64-41=", 2) | Out-Null
$d.Content.Find.Execute("51-43=", $true, $true, $false, $false, $false, $true, 1, $false, "4+21=", 2) | Out-Null
$d.Content.Find.Execute("11+4=", $true, $true, $false, $false, $false, $true, 1, $false, "84-17=", 2) | Out-Null
$d.Content.Find.Execute("22-9=", $true, $true, $false, $false, $false, $true, 1, $false, "38+34=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $true, $false, $false, $false, $true, 1, $false, "9+66=", 2) | Out-Null
$d.Content.Find.Execute("75-45=", $true, $true, $false, $false, $false, $true, 1, $false, "72-16=", 2) | Out-Null
$d.Content.Find.Execute("19+29=", $true, $true, $false, $false, $false, $true, 1, $false, "15+50=", 2) | Out-Null
$d.Content.Find.Execute("86-81=", $true, $true, $false, $false, $false, $true, 1, $false, "20+40=", 2) | Out-Null
$d.Content.Find.Execute("41-5=", $true, $true, $false, $false, $false, $true, 1, $false, "65+10=", 2) | Out-Null
$d.Content.Find.Execute("38+60=", $true, $true, $false, $false, $false, $true, 1, $false, "94-58=", 2) | Out-Null
$d.Content.Find.Execute("29-28=", $true, $true, $false, $false, $false, $true, 1, $false, "54-35=", 2) | Out-Null
$d.Content.Find.Execute("46-31=", $true, $true, $false, $false, $false, $true, 1, $false, "45-21=", 2) | Out-Null
$d.Content.Find.Execute("35+52=", $true, $true, $false, $false, $false, $true, 1, $false, "63-46=", 2) | Out-Null
$d.Content.Find.Execute("97-95=", $true, $true, $false, $false, $false, $true, 1, $false, "71-42=", 2) | Out-Null
$d.Content.Find.Execute("73-15=", $true, $true, $false, $false, $false, $true, 1, $false, "40+34=", 2) | Out-Null
$d.Content.Find.Execute("17+66=", $true, $true, $false, $false, $false, $true, 1, $false, "43-27=", 2) | Out-Null
